# Add Sheet2 with batch_size / accuracy data, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# Add a new worksheet after Sheet1 and rename it to "Sheet2"
$sheet1 = $wb.Worksheets.Item("Sheet1")
$newSheet = $wb.Worksheets.Add($null, $sheet1)
$newSheet.Name = "Sheet2"

$ws2 = $wb.Worksheets.Item("Sheet2")

# Header row
$ws2.Range("A1").Value = "batch_size"
$ws2.Range("B1").Value = "MNIST (uniform) test accuracy (%)"
$ws2.Range("C1").Value = "MNIST (normal) test accuracy (%)"
$ws2.Range("A1:C1").Font.Bold = $true

# Data rows
$data = @(
    @(16, 74.08, 68.25),
    @(32, 70.42, 67.24),
    @(64, 70.33, 69.28),
    @(128, 72.79, 74.61),
    @(256, 67.98, 72.13),
    @(512, 74.74, 67.59)
)

$row = 2
foreach ($r in $data) {
    $ws2.Cells.Item($row, 1).Value = $r[0]
    $ws2.Cells.Item($row, 2).Value = $r[1]
    $ws2.Cells.Item($row, 3).Value = $r[2]
    $row++
}

# Column widths similar to bestFit
$ws2.Columns.Item(1).ColumnWidth = 9.77734375
$ws2.Columns.Item(2).ColumnWidth = 30.44140625
$ws2.Columns.Item(3).ColumnWidth = 29.6640625

# Select C7 on Sheet2 and activate Sheet2 as the active tab
$ws2.Range("C7").Select() | Out-Null
$ws2.Activate() | Out-Null
